# The commit swaps the colour scheme carried by the presentation's theme
# part (ppt/theme/theme1.xml) from the "Integral / Red Violet" palette to
# the stock "Office Theme" palette (the palette that used to live only in
# ppt/theme/theme2.xml, the notes-master's theme part).
#
# PowerPoint's object model edits theme colours through
# ThemeColorScheme.Colors(index).RGB, where index 1-12 corresponds to the
# <a:clrScheme> child order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# RGB() itself isn't exposed as a builtin here, so reproduce VBA's packing
# (R + G*256 + B*65536) locally.

function ColorRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# New "Office Theme" palette (was previously only in ppt/theme/theme2.xml).
$colors.Colors(1).RGB  = ColorRGB 0x00 0x00 0x00   # dk1
$colors.Colors(2).RGB  = ColorRGB 0xFF 0xFF 0xFF   # lt1
$colors.Colors(3).RGB  = ColorRGB 0x44 0x54 0x6A   # dk2
$colors.Colors(4).RGB  = ColorRGB 0xE7 0xE6 0xE6   # lt2
$colors.Colors(5).RGB  = ColorRGB 0x5B 0x9B 0xD5   # accent1
$colors.Colors(6).RGB  = ColorRGB 0xED 0x7D 0x31   # accent2
$colors.Colors(7).RGB  = ColorRGB 0xA5 0xA5 0xA5   # accent3
$colors.Colors(8).RGB  = ColorRGB 0xFF 0xC0 0x00   # accent4
$colors.Colors(9).RGB  = ColorRGB 0x44 0x72 0xC4   # accent5
$colors.Colors(10).RGB = ColorRGB 0x70 0xAD 0x47   # accent6
$colors.Colors(11).RGB = ColorRGB 0x05 0x63 0xC1   # hlink
$colors.Colors(12).RGB = ColorRGB 0x95 0x4F 0x72   # folHlink
